$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Birds")
$ws1.Activate()
$ws1.Range("A271").Select()
$excel.ActiveWindow.FreezePanes = $true
$p2 = $excel.ActiveWindow.Panes.Item(2)
Write-Host ($p2 | Get-Member | Out-String)
